$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 61455
$ws.Range("B3").Value = 131577
$ws.Range("B4").Value = 123692
$ws.Range("B5").Value = 13390
$ws.Range("B6").Value = 28201
$ws.Range("B8").Value = 3118134
$ws.Range("B9").Value = 218681
$ws.Range("B10").Value = 29906
$ws.Range("B11").Value = 628817
$ws.Range("B12").Value = 325407
$ws.Range("B13").Value = 10773
$ws.Range("B14").Value = 186403
$ws.Range("B15").Value = 770842
$ws.Range("B16").Value = 3931
$ws.Range("B17").Value = 366305
$ws.Range("B18").Value = 1010987
$ws.Range("B19").Value = 12686
$ws.Range("B21").Value = 1190
$ws.Range("B22").Value = 314190
$ws.Range("B23").Value = 200693
$ws.Range("B24").Value = 48417
$ws.Range("B25").Value = 15082449
$ws.Range("B26").Value = 229
$ws.Range("B27").Value = 409495
$ws.Range("B28").Value = 13353
$ws.Range("B29").Value = 142903
$ws.Range("B30").Value = 4132
$ws.Range("B31").Value = 25837
$ws.Range("B32").Value = 18179
$ws.Range("B33").Value = 74946
$ws.Range("B34").Value = 1281846
$ws.Range("B35").Value = 6674
$ws.Range("B36").Value = 4862
$ws.Range("B37").Value = 1235778
$ws.Range("B38").Value = 102596
$ws.Range("B39").Value = 2968626
$ws.Range("B40").Value = 3854
$ws.Range("B41").Value = 11147
$ws.Range("B42").Value = 30240
$ws.Range("B43").Value = 265486
$ws.Range("B44").Value = 46344
$ws.Range("B45").Value = 342487
$ws.Range("B46").Value = 113876
$ws.Range("B47").Value = 68766
$ws.Range("B48").Value = 1642696
$ws.Range("B49").Value = 258204
$ws.Range("B50").Value = 11319
$ws.Range("B52").Value = 269958
$ws.Range("B53").Value = 396888
$ws.Range("B54").Value = 235140
$ws.Range("B55").Value = 70255
$ws.Range("B57").Value = 3742
$ws.Range("B58").Value = 124529
$ws.Range("B59").Value = 18474
$ws.Range("B60").Value = 261580
$ws.Range("B61").Value = 129
$ws.Range("B62").Value = 88332
$ws.Range("B63").Value = 5808421
$ws.Range("B64").Value = 23311
$ws.Range("B65").Value = 5925
$ws.Range("B66").Value = 319266
$ws.Range("B67").Value = 3517762
$ws.Range("B68").Value = 92856
$ws.Range("B69").Value = 358116
$ws.Range("B71").Value = 233696
$ws.Range("B72").Value = 22553
$ws.Range("B73").Value = 3739
$ws.Range("B74").Value = 14073
$ws.Range("B75").Value = 13164
$ws.Range("B76").Value = 218330
$ws.Range("B77").Value = 789188
$ws.Range("B78").Value = 6506
$ws.Range("B79").Value = 21892676
$ws.Range("B80").Value = 1703632
$ws.Range("B81").Value = 2627094
$ws.Range("B82").Value = 1103950
$ws.Range("B83").Value = 251904
$ws.Range("B84").Value = 838858
$ws.Range("B85").Value = 4092747
$ws.Range("B86").Value = 46428
$ws.Range("B87").Value = 628319
$ws.Range("B88").Value = 718632
$ws.Range("B89").Value = 393639
$ws.Range("B90").Value = 162666
$ws.Range("B91").Value = 126745
$ws.Range("B92").Value = 105784
$ws.Range("B93").Value = 282981
$ws.Range("B94").Value = 97663
$ws.Range("B95").Value = 1205
$ws.Range("B96").Value = 122545
$ws.Range("B97").Value = 531834
$ws.Range("B98").Value = 10761
$ws.Range("B99").Value = 2114
$ws.Range("B100").Value = 179697
$ws.Range("B101").Value = 2963
$ws.Range("B102").Value = 255384
$ws.Range("B103").Value = 68153
$ws.Range("B104").Value = 38641
$ws.Range("B105").Value = 34158
$ws.Range("B106").Value = 432425
$ws.Range("B107").Value = 34134
$ws.Range("B108").Value = 14059
$ws.Range("B109").Value = 30432
$ws.Range("B111").Value = 18613
$ws.Range("B112").Value = 1226
$ws.Range("B113").Value = 2361874
$ws.Range("B114").Value = 252413
$ws.Range("B115").Value = 2477
$ws.Range("B116").Value = 44016
$ws.Range("B117").Value = 98142
$ws.Range("B118").Value = 513314
$ws.Range("B119").Value = 70138
$ws.Range("B120").Value = 49552
$ws.Range("B121").Value = 377603
$ws.Range("B122").Value = 1572861
$ws.Range("B123").Value = 2634
$ws.Range("B125").Value = 5310
$ws.Range("B126").Value = 165340
$ws.Range("B127").Value = 115818
$ws.Range("B128").Value = 199344
$ws.Range("B129").Value = 854240
$ws.Range("B130").Value = 366762
$ws.Range("B131").Value = 11630
$ws.Range("B132").Value = 294233
$ws.Range("B133").Value = 1824457
$ws.Range("B134").Value = 1087885
$ws.Range("B135").Value = 2824425
$ws.Range("B136").Value = 838852
$ws.Range("B137").Value = 210070
$ws.Range("B138").Value = 1063949
$ws.Range("B139").Value = 4808133
$ws.Range("B140").Value = 25539
$ws.Range("B142").Value = 4607
$ws.Range("B143").Value = 1895
$ws.Range("B146").Value = 2317
$ws.Range("B147").Value = 424445
$ws.Range("B148").Value = 40621
$ws.Range("B149").Value = 698518
$ws.Range("B152").Value = 61311
$ws.Range("B153").Value = 385022
$ws.Range("B154").Value = 245159
$ws.Range("B156").Value = 14368
$ws.Range("B157").Value = 1592626
$ws.Range("B158").Value = 10637
$ws.Range("B159").Value = 3567408
$ws.Range("B160").Value = 121338
$ws.Range("B161").Value = 34461
$ws.Range("B162").Value = 10933
$ws.Range("B163").Value = 1007792
$ws.Range("B164").Value = 670613
$ws.Range("B165").Value = 23256
$ws.Range("B166").Value = 1178
$ws.Range("B169").Value = 78855
$ws.Range("B170").Value = 2965
$ws.Range("B171").Value = 13106
$ws.Range("B172").Value = 12720
$ws.Range("B173").Value = 318236
$ws.Range("B174").Value = 4998089
$ws.Range("B175").Value = 32651865
$ws.Range("B176").Value = 42224
$ws.Range("B177").Value = 2160809
$ws.Range("B178").Value = 532710
$ws.Range("B179").Value = 4446752
$ws.Range("B180").Value = 216146
$ws.Range("B181").Value = 93597
$ws.Range("B183").Value = 205181
$ws.Range("B184").Value = 3137
$ws.Range("B185").Value = 300946
$ws.Range("B186").Value = 6446
$ws.Range("B187").Value = 92004
$ws.Range("B188").Value = 38403
